$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.232.72"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "1.662.73"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").Value = "218.36"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "0.5219"
$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("D8").Value = "0.2654"
$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").Value = "0.06286"
$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("D10").Value = "20.82"
$ws.Range("E10").Value = "  -4.58%  "

$ws.Range("D11").Value = "0.07729"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").Value = "1.667.39"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").Value = "4.431"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").Value = "1.891.53"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").Value = "0.5444"
$ws.Range("E15").Value = "  -2.44%  "

$ws.Range("D16").Value = "0.0₅8198"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("D17").Value = "64.58"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").Value = "26.281.30"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D20").Value = "4.652"
$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("D21").Value = "193.46"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  -2.04%  "

$ws.Range("D23").Value = "6.038"
$ws.Range("E23").Value = "  -4.86%  "

$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("D25").Value = "139.82"
$ws.Range("E25").Value = "  -1.93%  "

$ws.Range("D26").Value = "0.1229"
$ws.Range("E26").Value = "  -4.67%  "

$ws.Range("D27").Value = "7.157"
$ws.Range("E27").Value = "  -3.24%  "

$ws.Range("D28").Value = "16.10"
$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").Value = "1.416"
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("D30").Value = "0.06118"
$ws.Range("E30").Value = "  -3.15%  "

$ws.Range("D31").Value = "1.281"
$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").Value = "3.571"
$ws.Range("E32").Value = "  -1.81%  "

$ws.Range("D33").Value = "3.254"
$ws.Range("E33").Value = "  -5.81%  "

$ws.Range("D34").Value = "1.620"
$ws.Range("E34").Value = "  -3.54%  "

$ws.Range("D35").Value = "0.9641"
$ws.Range("E35").Value = "  -4.57%  "

$ws.Range("D36").Value = "2.424"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "2.789"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Value = "0.5713"
$ws.Range("E38").Value = "  -7.90%  "

$ws.Range("D39").Value = "0.01607"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").Value = "5.990"
$ws.Range("E40").Value = "  -2.83%  "

$ws.Range("D41").Value = "0.8566"

$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("D43").Value = "1.015.02"
$ws.Range("E43").Value = "  -6.65%  "

$ws.Range("D44").Value = "100.08"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").Value = "1.807.31"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  +7.92%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "57.14"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("D49").Value = "8.008"
$ws.Range("E49").Value = "  -2.27%  "

$ws.Range("D50").Value = "1.480"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").Value = "0.05185"
$ws.Range("E51").Value = "  -0.43%  "
